$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on column D so numeric-looking price strings (e.g. "1.000",
# "245.00") are preserved exactly as literal text instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.492.75'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.727.10'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '245.00'
$ws.Range('E5').Value = '  +2.33%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4802'
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('D8').Value = '0.2673'
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('D9').Value = '0.06215'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = '1.727.08'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').Value = '0.07152'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('D13').Value = '0.6178'
$ws.Range('E13').Value = '  +4.58%  '
$ws.Range('D14').Value = '4.509'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '26.507.73'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '0.000006941'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').Value = '11.65'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').Value = '1.949.04'
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').Value = '4.529'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '8.950'
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('D24').Value = '5.284'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').Value = '136.54'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').Value = '1.409'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').Value = '106.71'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '3.974'
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('E31').Value = '  +3.84%  '
$ws.Range('D32').Value = '3.713'
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').Value = '0.04562'
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').Value = '0.9995'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.615'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6355'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9886'
$ws.Range('E37').Value = '  +1.71%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '0.9357'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.089'
$ws.Range('E39').Value = '  +9.59%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.416'
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = '105.12'
$ws.Range('E41').Value = '  -8.73%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01502'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.672'
$ws.Range('E44').Value = '  +7.30%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3905'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.895'
$ws.Range('E46').Value = '  +10.49%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1189'
$ws.Range('E47').Value = '  +3.48%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05329'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '31.01'
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.882'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.265'
$ws.Range('E51').Value = '  +3.69%  '
